# Removes the two "Night School stop..." and "Define constants for
# amounts..." to-do items from the list, per the commit's "Updated ToDo
# list". The trailing _GoBack bookmark (bookmarkStart/bookmarkEnd) that
# used to sit at the end of the "Define constants..." paragraph ends up
# attached to the end of the surviving "...Excel sheet" paragraph, since
# the paragraph marks in between are what actually get deleted.

$d = $word.ActiveDocument

function Find-ParaIndex($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "$prefix*") {
            return $i
        }
    }
    return -1
}

# 1) Delete the entire "Night School stop: ... ??" paragraph, including its
#    own paragraph mark, so the list simply closes the gap.
$nightIdx = Find-ParaIndex("Night School stop:")
$pNight = $d.Paragraphs.Item($nightIdx)
$d.Range($pNight.Range.Start, $pNight.Range.End).Delete()

# 2) Clear the text of the "Define constants for amounts..." paragraph but
#    keep its paragraph mark — that mark is what still carries the
#    _GoBack bookmark, and it must survive.
$defineIdx = Find-ParaIndex("Define constants for amounts")
$pDefine = $d.Paragraphs.Item($defineIdx)
$d.Range($pDefine.Range.Start, $pDefine.Range.End - 1).Delete()

# 3) Merge that now-empty (but bookmark-bearing) paragraph mark up into the
#    preceding "...Excel sheet" paragraph by deleting the Excel-sheet
#    paragraph's own trailing mark.
$excelIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Excel sheet*") {
        $excelIdx = $i
    }
}
$pExcel = $d.Paragraphs.Item($excelIdx)
$d.Range($pExcel.Range.End - 1, $pExcel.Range.End).Delete()
